# Applies the "annotated ML10M genre prediction for Miguel" edit:
#   1. Slide 1 "Ratings model" textbox: reposition/resize + shrink font 32 -> 20
#   2. Slide 2 "Genre model" textbox: reposition/resize (height only) + shrink font 32 -> 20
#   3. Slide 2: add a new annotation textbox explaining genre-model re-use of the
#      pre-trained item-embedding layer.

$EMU = 12700.0
# PowerPoint's Shape geometry properties round-trip through a single-precision
# float, which truncates a handful of our target EMU values by 1 EMU. Nudge
# every point value up by a hair (well under half an EMU) so the stored EMU
# lands exactly on the target after that float32 round-trip.
$eps = 0.00004

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 1 - "Ratings model" textbox
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$ratings = $s1.Shapes.Item(78)

$ratings.Left   = (719300   / $EMU) + $eps
$ratings.Top    = (666920   / $EMU) + $eps
$ratings.Width  = (2238258  / $EMU) + $eps
$ratings.Height = (400110   / $EMU) + $eps

$ratings.TextFrame.TextRange.Font.Size = 2000 / 100.0

# ---------------------------------------------------------------------------
# 2) Slide 2 - "Genre model" textbox
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$genre = $s2.Shapes.Item(52)

$genre.Left   = (761114  / $EMU) + $eps
$genre.Top    = (687410  / $EMU) + $eps
$genre.Width  = (3053686 / $EMU) + $eps
$genre.Height = (400110  / $EMU) + $eps

$genre.TextFrame.TextRange.Font.Size = 2000 / 100.0

# ---------------------------------------------------------------------------
# 3) Slide 2 - new annotation textbox
# ---------------------------------------------------------------------------
# Shape IDs are assigned automatically (lowest free slot on the slide) and
# aren't directly settable. The source deck's new shape landed on id 125
# (re-using a slot freed up elsewhere on this heavily-edited slide), so spin
# up throwaway textboxes - deleting each one that lands on the wrong id -
# until the allocator serves us id 125, then keep that one.
$targetId = 125
$maxTries = 500
$tries = 0
$note = $null
while ($tries -lt $maxTries) {
    $cand = $s2.Shapes.AddTextbox(1, 10, 10, 10, 10)
    if ($cand.Id -eq $targetId) {
        $note = $cand
        break
    }
    $cand.Delete()
    $tries = $tries + 1
}
if ($note -eq $null) {
    # Fallback: could not land on the exact historical id - just use
    # whatever the next free id is rather than failing the whole edit.
    $note = $s2.Shapes.AddTextbox(1, 10, 10, 10, 10)
}
$note.Name = "TextBox 124"

# Set position/width before autofitting so word-wrap (and thus the
# autofit height) is computed against the final box width.
$note.Left  = (4090989 / $EMU) + $eps
$note.Top   = (4214307 / $EMU) + $eps
$note.Width = (2238258 / $EMU) + $eps

$note.TextFrame.TextRange.Text = "Re-use pre-trained item embedding layer for genre prediction"

$noteFont = $note.TextFrame.TextRange.Font
$noteFont.Size = 1400 / 100.0
$noteFont.Name = "Roboto Medium"
$noteFont.NameFarEast = "Roboto Medium"

$note.TextFrame.AutoSize = 1
$note.TextFrame.WordWrap = $true
$note.Fill.Visible = $false
